# Auto-generated edit script
# Part 1: fix delimiter in "Razon social"/"Nombre Fantasia" text fields
#   - replace comma separators between co-contractors with periods
#   - normalize "S.H." abbreviation to "SH"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E93").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E95").Value = "MONTICO. RICARDO"
$ws.Range("E129").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E135").Value = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
$ws.Range("F143").Value = "MORERA. SALVADOR"
$ws.Range("E189").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E201").Value = "OLIVERA. FLORENCIO"
$ws.Range("F201").Value = "OLIVERA. FLORENCIO"
$ws.Range("E204").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("F219").Value = "PARRAVICINI VIRGINIA VANINA. VIRGINIA VANINA"
$ws.Range("E224").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E294").Value = "ALVAREZ. RENZO JOEL"
$ws.Range("F294").Value = "ALVAREZ. RENZO JOEL"
$ws.Range("E296").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# Part 2: fix floating point formatting in "Importe" column (H)
#   European format "1.234,56" (thousands-dot, decimal-comma)
#   ->  plain decimal "1234.56" (no thousands separator, dot decimal)
#   Cells keep their original text (string) storage, so the column
#   number format is switched to Text ("@") before assignment -
#   otherwise Excel would auto-coerce the text into a Double and
#   silently drop the trailing zeros / fixed decimal formatting.
$importeFixes = @(
    @(2, "239000.00"),
    @(3, "59199.36"),
    @(4, "440000.00"),
    @(5, "5599.96"),
    @(6, "2250.00"),
    @(7, "1760.00"),
    @(8, "17780.00"),
    @(9, "57620.00"),
    @(10, "2688.36"),
    @(11, "21600.00"),
    @(12, "2320.00"),
    @(13, "3800.00"),
    @(14, "1121530.43"),
    @(15, "62280.00"),
    @(16, "131190.00"),
    @(17, "58779.86"),
    @(18, "1544344.46"),
    @(19, "8934.00"),
    @(20, "169640.00"),
    @(21, "16341.70"),
    @(22, "43509.18"),
    @(23, "201413.80"),
    @(24, "22500.00"),
    @(25, "288.00"),
    @(26, "472.00"),
    @(27, "5775.00"),
    @(28, "19.80"),
    @(29, "848.00"),
    @(30, "116.85"),
    @(31, "5750.00"),
    @(32, "31.40"),
    @(33, "3018533.03"),
    @(34, "113231.64"),
    @(35, "24349.00"),
    @(36, "311968.25"),
    @(37, "6952.00"),
    @(38, "6520.00"),
    @(39, "6565.70"),
    @(40, "19950.00"),
    @(41, "13040.00"),
    @(42, "1000.00"),
    @(43, "51051.09"),
    @(44, "65005.57"),
    @(45, "1681.58"),
    @(46, "19250.00"),
    @(47, "15230.38"),
    @(48, "390.00"),
    @(49, "32758.53"),
    @(50, "8926.92"),
    @(51, "60457.05"),
    @(52, "7923.30"),
    @(53, "940.00"),
    @(54, "2032.01"),
    @(55, "900.00"),
    @(56, "19949.00"),
    @(57, "1098.00"),
    @(58, "240.00"),
    @(59, "478488.60"),
    @(60, "19980.00"),
    @(61, "1079458.25"),
    @(62, "34972.55"),
    @(63, "7339.28"),
    @(64, "25175.69"),
    @(65, "150.00"),
    @(66, "33863.80"),
    @(67, "700.00"),
    @(68, "47577.36"),
    @(69, "6930.00"),
    @(70, "16700.00"),
    @(71, "87321.47"),
    @(72, "19687.50"),
    @(73, "13020.00"),
    @(74, "1771.90"),
    @(75, "4380.00"),
    @(76, "4942.00"),
    @(77, "32476.00"),
    @(78, "4240.00"),
    @(79, "23000.00"),
    @(80, "20580.00"),
    @(81, "1500.00"),
    @(82, "195000.00"),
    @(83, "3150.00"),
    @(84, "24000.00"),
    @(85, "800.00"),
    @(86, "13000.00"),
    @(87, "17075.00"),
    @(88, "124080.00"),
    @(89, "11900.00"),
    @(90, "890.00"),
    @(91, "57960.00"),
    @(92, "38432.16"),
    @(93, "3990.00"),
    @(94, "2200.00"),
    @(95, "11900.00"),
    @(96, "9590.00"),
    @(97, "36072.65"),
    @(98, "9980.00"),
    @(99, "9150.00"),
    @(100, "56099.98"),
    @(101, "5219.50"),
    @(102, "1300.00"),
    @(103, "12374.86"),
    @(104, "6910.00"),
    @(105, "1239.75"),
    @(106, "10500.00"),
    @(107, "29380.00"),
    @(108, "102835.56"),
    @(109, "36958.45"),
    @(110, "1650.00"),
    @(111, "1980.00"),
    @(112, "4992.98"),
    @(113, "7322.47"),
    @(114, "1470.00"),
    @(115, "90.00"),
    @(116, "6335.00"),
    @(117, "12524.00"),
    @(118, "46980.00"),
    @(119, "10160.00"),
    @(120, "6716.00"),
    @(121, "4410.96"),
    @(122, "51886.60"),
    @(123, "20072.50"),
    @(124, "28369.00"),
    @(125, "10542.05"),
    @(126, "6100.00"),
    @(127, "47916.50"),
    @(128, "35920.00"),
    @(129, "23030.00"),
    @(130, "13500.00"),
    @(131, "66550.00"),
    @(132, "75000.00"),
    @(133, "391.20"),
    @(134, "5185.00"),
    @(135, "17800.00"),
    @(136, "850.00"),
    @(137, "91634.00"),
    @(138, "4620.00"),
    @(139, "39460.00"),
    @(140, "8000.00"),
    @(141, "10000.00"),
    @(142, "5000.00"),
    @(143, "2000.00"),
    @(144, "5000.00"),
    @(145, "30000.00"),
    @(146, "29000.00"),
    @(147, "10000.00"),
    @(148, "84938.66"),
    @(149, "10866.00"),
    @(150, "2629.74"),
    @(151, "171676.83"),
    @(152, "5274.32"),
    @(153, "43407.00"),
    @(154, "26369.80"),
    @(155, "43410.00"),
    @(156, "16500.00"),
    @(157, "22000.00"),
    @(158, "22000.00"),
    @(159, "20000.00"),
    @(160, "7000.00"),
    @(161, "12000.00"),
    @(162, "4000.00"),
    @(163, "31000.00"),
    @(164, "10000.00"),
    @(165, "10000.00"),
    @(166, "10000.00"),
    @(167, "9000.00"),
    @(168, "10500.00"),
    @(169, "8000.00"),
    @(170, "10000.00"),
    @(171, "6000.00"),
    @(172, "30000.00"),
    @(173, "10000.00"),
    @(174, "23500.00"),
    @(175, "10000.00"),
    @(176, "10000.00"),
    @(177, "20870.00"),
    @(178, "1000.00"),
    @(179, "5000.00"),
    @(180, "75166.50"),
    @(181, "18000.00"),
    @(182, "10000.00"),
    @(183, "12000.00"),
    @(184, "12000.00"),
    @(185, "14000.00"),
    @(186, "23205.00"),
    @(187, "40000.00"),
    @(188, "19000.00"),
    @(189, "10000.00"),
    @(190, "5000.00"),
    @(191, "10000.00"),
    @(192, "10000.00"),
    @(193, "78000.00"),
    @(194, "51600.00"),
    @(195, "45000.00"),
    @(196, "23150.00"),
    @(197, "69700.00"),
    @(198, "1500.00"),
    @(199, "19120.00"),
    @(200, "200108.00"),
    @(201, "22000.00"),
    @(202, "9870.00"),
    @(203, "69.00"),
    @(204, "18705.00"),
    @(205, "95800.00"),
    @(206, "1990.00"),
    @(207, "71140.00"),
    @(208, "66740.00"),
    @(209, "8971.00"),
    @(210, "1040.70"),
    @(211, "6150.00"),
    @(212, "38630.00"),
    @(213, "112900.00"),
    @(214, "4500.00"),
    @(215, "254723.00"),
    @(216, "30160.00"),
    @(217, "7980.00"),
    @(218, "14450.00"),
    @(219, "35226.21"),
    @(220, "13346.97"),
    @(221, "36600.00"),
    @(222, "56400.00"),
    @(223, "18005.90"),
    @(224, "49290.00"),
    @(225, "4980.00"),
    @(226, "1300.00"),
    @(227, "4794.21"),
    @(228, "10800.00"),
    @(229, "862.56"),
    @(230, "8259.37"),
    @(231, "30487.02"),
    @(232, "19200.00"),
    @(233, "19540.00"),
    @(234, "1000.00"),
    @(235, "17600.00"),
    @(236, "13202.60"),
    @(237, "29550.00"),
    @(238, "1537.78"),
    @(239, "317009.00"),
    @(240, "80000.00"),
    @(241, "40000.00"),
    @(242, "40000.00"),
    @(243, "40000.00"),
    @(244, "80000.00"),
    @(245, "40000.00"),
    @(246, "40000.00"),
    @(247, "40000.00"),
    @(248, "40000.00"),
    @(249, "80000.00"),
    @(250, "80000.00"),
    @(251, "8500.00"),
    @(252, "89684.50"),
    @(253, "9209996.57"),
    @(254, "4000.00"),
    @(255, "600.00"),
    @(256, "18600.00"),
    @(257, "24800.00"),
    @(258, "12870372.93"),
    @(259, "2436350.88"),
    @(260, "332000.00"),
    @(261, "377520.00"),
    @(262, "315720.00"),
    @(263, "332000.00"),
    @(264, "309000.00"),
    @(265, "309000.00"),
    @(266, "588000.00"),
    @(267, "309000.00"),
    @(268, "807350.00"),
    @(269, "1046000.00"),
    @(270, "350720.00"),
    @(271, "309000.00"),
    @(272, "309000.00"),
    @(273, "618000.00"),
    @(274, "644100.00"),
    @(275, "610400.00"),
    @(276, "910400.00"),
    @(277, "588000.00"),
    @(278, "908980.00"),
    @(279, "618000.00"),
    @(280, "339143.00"),
    @(281, "397225.00"),
    @(282, "480816.63"),
    @(283, "130000.00"),
    @(284, "18501572.00"),
    @(285, "16886584.01"),
    @(286, "484000.00"),
    @(287, "14000.00"),
    @(288, "500.00"),
    @(289, "1100.00"),
    @(290, "3117166.68"),
    @(291, "18900.00"),
    @(292, "8000.00"),
    @(293, "1734.53"),
    @(294, "23700.00"),
    @(295, "393220.00"),
    @(296, "19000.00"),
    @(297, "72502.00"),
    @(298, "10650.00"),
    @(299, "600.00")
)

$importeCol = 8  # column H
foreach ($fix in $importeFixes) {
    $r = $fix[0]
    $newValue = $fix[1]
    $cell = $ws.Cells.Item($r, $importeCol)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}

Write-Output "Updated $($importeFixes.Count) Importe cells and 14 Razon social/Nombre Fantasia cells."
